# Services.xlsx – clear the stray "sales price" (column D) values that were
# mistakenly duplicated from the "price" column (column C) for the three
# "Free"-priced rows: row 5 (hvac-install Proposal), row 7 (Other / contact),
# and row 8 (Vinal-Flooring-Proposal). The "price" column keeps "Free"; only
# "sales price" is cleared back to blank, leaving the cell style intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("D8").ClearContents()

# Leave the sheet scrolled to / selected on the last cell the author touched
# (D7) so the saved view matches the final state of the workbook.
$ws.Range("D7").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
